# Aggregated Monthly Infographic - month roll-forward (Sep -> Oct 2024)
# and label rename: "Ordering Experience" -> "Customer Service" -> "Food"
# (Save as PDF and Email functionality)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the LABEL lookup text (column E refers to these shared strings) ---
$ws.Range("E27").Value = "Customer Service"
$ws.Range("E28").Value = "Food"

# --- Per-row data refresh: MO (col C) 9 -> 10, and STATISTIC (col D) values ---
$rowData = @(
    @(2, "125"),
    @(3, "121"),
    @(4, "96.42"),
    @(5, "7"),
    @(6, "0.06"),
    @(7, "378"),
    @(8, "2.03"),
    @(9, "45"),
    @(10, "2.71"),
    @(11, "8"),
    @(12, "4"),
    @(13, "16"),
    @(14, "40"),
    @(15, "14"),
    @(16, "21"),
    @(17, "53"),
    @(18, "4,621"),
    @(19, "1,883"),
    @(20, "493"),
    @(21, "1,011"),
    @(22, "379"),
    @(23, "0.261816"),
    @(24, "0.536909"),
    @(25, "0.201275"),
    @(26, "432"),
    @(27, "394"),
    @(28, "379"),
    @(29, "648"),
    @(30, "0.85"),
    @(31, "0"),
    @(32, "0"),
    @(33, "0"),
    @(34, "0"),
    @(35, "0"),
    @(36, "0"),
    @(37, "0"),
    @(38, "0"),
    @(39, "0"),
    @(40, "0"),
    @(41, "0"),
    @(42, "0"),
    @(43, "0"),
    @(44, "0"),
    @(45, "0"),
    @(46, "0"),
    @(47, "0"),
    @(48, "0"),
    @(49, "0"),
    @(50, "0"),
    @(51, "0"),
    @(52, "0"),
    @(53, "0"),
    @(54, "0"),
    @(55, "0")
)

foreach ($item in $rowData) {
    $r = $item[0]
    $newVal = $item[1]

    # Month column: September (9) -> October (10)
    $ws.Cells.Item($r, 3).Value = 10

    # Statistic column: refresh the quoted-text formula with the new figure,
    # and drop the inherited number-format style (now unformatted/general).
    $cell = $ws.Cells.Item($r, 4)
    $cell.Style = "Normal"
    $cell.Formula = '="' + $newVal + '"'
}
